# Apply "US 3.3 commit files" changes to
# InputData/bldgs/PPEIdtIL/Potential Perc Eff Improvement due to Impr Labeling.xlsx

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("PPEIdtIL")

# ---------------------------------------------------------------------------
# "About" sheet: insert a new explanatory note (rows 10-16), pushing the
# pre-existing ACEEE source note down to rows 18-22 (row 17 left blank as a
# separator, row 9 "Note:" header untouched).
# ---------------------------------------------------------------------------

$oldNote1 = $wsAbout.Range("A10").Value2
$oldNote2 = $wsAbout.Range("A11").Value2
$oldNote3 = $wsAbout.Range("A12").Value2
$oldNote4 = $wsAbout.Range("A13").Value2
$oldNote5 = $wsAbout.Range("A14").Value2

$wsAbout.Range("A18").Value2 = $oldNote1
$wsAbout.Range("A19").Value2 = $oldNote2
$wsAbout.Range("A20").Value2 = $oldNote3
$wsAbout.Range("A21").Value2 = $oldNote4
$wsAbout.Range("A22").Value2 = $oldNote5

$wsAbout.Range("A10").Value2 = "This variable reflects improvement in efficiency components selected by consumers due"
$wsAbout.Range("A11").Value2 = "to improved labeling. The labeling influences consumers who are buying appliances of all"
$wsAbout.Range("A12").Value2 = "Quality levels, so it's represented as a simple percentage increase in the efficiency of"
$wsAbout.Range("A13").Value2 = "components sold (at all quality levels). If Quality Levels are defined based on"
$wsAbout.Range("A14").Value2 = "particular efficiency thresholds, this may mean the number of square feet served by"
$wsAbout.Range("A15").Value2 = "components of a given quality level will not be accurate. It's just a question of the meaning"
$wsAbout.Range("A16").Value2 = "of the labels given to each Quality Level."

# ---------------------------------------------------------------------------
# "PPEIdtIL" sheet: header text + wrap formatting, and updated percentages.
# ---------------------------------------------------------------------------

$wsData.Range("A1").Value2 = "Efficiency Improvement by Building Component (dimensionless)"
$wsData.Range("A1").WrapText = $true
$wsData.Rows.Item(1).RowHeight = 45

$wsData.Range("B2").Value2 = 0.02
$wsData.Range("D2").Value2 = 0.02

$wsData.Range("B3").Value2 = 0.02
$wsData.Range("D3").Value2 = 0.02

$wsData.Range("B6").Value2 = 0.02
$wsData.Range("D6").Value2 = 0.02

# ---------------------------------------------------------------------------
# Active sheet / selection bookkeeping: PPEIdtIL becomes the active tab with
# selection L3, while About keeps a (non-selected) selection of E35.
# ---------------------------------------------------------------------------

$wsAbout.Range("E35").Select() | Out-Null
$wsData.Activate()
$wsData.Range("L3").Select() | Out-Null
